$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Batch 1: rows 49-64.
# Cell assignment order below is chosen so that newly created shared-string
# entries come out in the same order as in the target workbook (indices 49-67).
$ws.Range("B49").Value = "Chức năng xóa"
$ws.Range("D50").Value = "phần 1 :"
$ws.Range("D55").Value = "Phần 2 :"
$ws.Range("F55").Value = "gán text xóa vào link"
$ws.Range("F56").Value = "nhacungcap_list.php?xoa=true&MaNhaCungCap=6"
$ws.Range("F50").Value = "kiểm tra có phải biến xóa = true"
$ws.Range("F51").Value = 1
$ws.Range("G51").Value = "và mã nhà cung cấp có giá trị k?"
$ws.Range("F52").Value = 2
$ws.Range("G52").Value = "Viết câu lệnh sql xóa theo mã nhà cung cấp"
$ws.Range("F53").Value = 3
$ws.Range("G53").Value = "thực hiện xoa"
$ws.Range("F54").Value = 4
$ws.Range("G54").Value = "Xóa thành công thì hiển thị thông báo"
$ws.Range("F57").Value = "add đoạn js xử lý : hỏi trước khi xóa"
$ws.Range("B58").Value = "chức năng tìm kiếm tên nhà cung cấp"
$ws.Range("D59").Value = "Phần HTML "
$ws.Range("E60").Value = "thêm 1 dòng"
$ws.Range("H60").Value = "chứa label, textbox, button"
$ws.Range("E61").Value = "và form"
$ws.Range("F62").Value = "action"
$ws.Range("H62").Value = "nhacungcap_list"
$ws.Range("F63").Value = "method"
$ws.Range("H63").Value = "GET"
$ws.Range("D64").Value = "Phần PHP"

# Commit this batch so row "spans" are recomputed together (matches target: 2:8).
$wb.Save()

# Batch 2: rows 65-67, committed separately (matches target spans: 5:6).
$ws.Range("E65").Value = "kiểm tra điều kiện search"
$ws.Range("F66").Value = "add điều kiện where sử dụng like "
$ws.Range("E67").Value = "nối chuỗi vào câu sql"

# Update the view/selection to match the target workbook state.
$ws.Range("D68").Select() | Out-Null
